# Regenerate the handoff report: the markdown source file was renamed
# (new GUID) and its handoff to zh-cn/de-de failed, so both locale sheets
# lose their "latest handoff" file/date and are now flagged Ignored with
# a reset handback timestamp.

$wb = $excel.ActiveWorkbook

$oldFile = "3629b8fc-7e47-4fb7-bc9a-4a2733432db0.md"
$newFile = "c5c317c1-99b9-4fdc-9fbb-b71ed502eb84.md"
$configFile = ".localization-config"
$defaultDate = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# Overview sheet: rename the source file and mark it as failed handoff.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("A2").Value = $newFile
$overview.Range("B2").Value = "Handoff failed"
$overview.Range("C2").Value = "Handoff failed"

$overview.Hyperlinks.Delete()
$overview.Hyperlinks.Add($overview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/9e818770b5e8f2a3ff3c153bf5635632227254d1/e2e/" + $newFile, "", "", $newFile) | Out-Null
$overview.Hyperlinks.Add($overview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/9e818770b5e8f2a3ff3c153bf5635632227254d1/" + $configFile, "", "", $configFile) | Out-Null

# ---------------------------------------------------------------------
# zh-cn sheet: handoff failed -> drop the handoff file/date, reset the
# handback date, and flag the row as Ignored.
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("A2").Value = $newFile
$zhcn.Range("B2").Value = "Handoff failed"
$zhcn.Range("C2").Clear()
$zhcn.Range("D2").Value = $defaultDate
$zhcn.Range("G2").Value = $defaultDate
$zhcn.Range("H2").Value = "Ignored"

$zhcn.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/9e818770b5e8f2a3ff3c153bf5635632227254d1/e2e/" + $newFile, "", "", $newFile) | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/9e818770b5e8f2a3ff3c153bf5635632227254d1/" + $configFile, "", "", $configFile) | Out-Null

# ---------------------------------------------------------------------
# de-de sheet: same treatment as zh-cn.
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("A2").Value = $newFile
$dede.Range("B2").Value = "Handoff failed"
$dede.Range("C2").Clear()
$dede.Range("D2").Value = $defaultDate
$dede.Range("G2").Value = $defaultDate
$dede.Range("H2").Value = "Ignored"

$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/9e818770b5e8f2a3ff3c153bf5635632227254d1/e2e/" + $newFile, "", "", $newFile) | Out-Null
$dede.Hyperlinks.Add($dede.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/9e818770b5e8f2a3ff3c153bf5635632227254d1/" + $configFile, "", "", $configFile) | Out-Null
